$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'267.51"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'22.82"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'6.330"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.06201"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.594"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.670"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'1.391"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8290"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.01362"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.1612"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.08248"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03413"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.03153"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.09289"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'3.909"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.001730"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.04839"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.006315"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.005386"
$ws.Range("D20").ClearFormats()
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = "'0.001091"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = "'3.761"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = "'2.369"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = "'0.3350"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = "'0.1214"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = "'0.0002685"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '26UpBotsUBXT'
$ws.Range("D40").Value = "'0.04654"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.006871"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.1156"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'0.003463"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.01220"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00006273"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.7896"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.1624"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.01241"
$ws.Range("D50").ClearFormats()
